# Simu Right turn 2,2g rigidité châssis
# Ajout des Simu right turn 2,2g pour la rigidité du châssis

$wb = $excel.ActiveWorkbook

# --- Fill in the "Right Turn 2,2G" sheet with the simulation results ---
$ws = $wb.Worksheets.Item("Right Turn 2,2G")

# Title cell (merged D2:E3)
$ws.Range("D2").Value = "Right Turn 2,2G"

# Data table D6:F15
$ws.Range("D6").Value = -0.605
$ws.Range("E6").Value = 637
$ws.Range("F6").Value = 5.27

$ws.Range("D7").Value = 0.461
$ws.Range("E7").Value = 6.37
$ws.Range("F7").Value = 5.84

$ws.Range("D8").Value = -0.586
$ws.Range("E8").Value = 5.45
$ws.Range("F8").Value = 4.33

$ws.Range("D9").Value = 0.474
$ws.Range("E9").Value = 5.44
$ws.Range("F9").Value = 4.78

$ws.Range("D10").Value = -0.588
$ws.Range("E10").Value = 4.39
$ws.Range("F10").Value = 3.01

$ws.Range("D11").Value = 0.445
$ws.Range("E11").Value = 4.41
$ws.Range("F11").Value = 3.58

$ws.Range("D12").Value = -0.419
$ws.Range("E12").Value = 2.15
$ws.Range("F12").Value = 0.916

$ws.Range("D13").Value = 0.636
$ws.Range("E13").Value = 2.13
$ws.Range("F13").Value = 1.55

$ws.Range("D14").Value = -0.392
$ws.Range("E14").Value = 1.42
$ws.Range("F14").Value = 0.035

$ws.Range("D15").Value = 0.382
$ws.Range("E15").Value = 1.38
$ws.Range("F15").Value = 0.903

# Selection on this sheet changes and it becomes the active/selected tab
$ws.Range("H13").Select() | Out-Null

# Make "Right Turn 2,2G" the active sheet (it becomes tabSelected and the
# workbook's active tab)
$ws.Activate() | Out-Null

# --- Update the selection on the "MAX SPEED" sheet as well ---
$wsMax = $wb.Worksheets.Item("MAX SPEED")
$wsMax.Range("E32").Select() | Out-Null

# Re-activate "Right Turn 2,2G" so it remains the workbook's active tab
$ws.Activate() | Out-Null

$wb.Save() | Out-Null
